# Yearly coverage in scenario 3b
#
# Fill in the previously-missing yearly coverage values (0.6) on the
# "Platform Coverage" sheet so that every year column from H2 to AD2
# (2018-2040) has a value, instead of only every other year. Also bring
# "Platform Coverage" to the front as the active sheet/selection, matching
# the saved view state.

$wb = $excel.ActiveWorkbook

$wsCoverage = $wb.Worksheets.Item("Platform Coverage")

# Fill every column from H to AD on row 2 with 0.6 (previously only every
# other column, H/J/L/.../AD, was populated).
$wsCoverage.Range("H2:AD2").Value = 0.6

# Make "Platform Coverage" the active / selected sheet, scrolled so that
# column R is left-most visible, with AE2 as the active selection.
$wsCoverage.Activate()
$excel.ActiveWindow.ScrollColumn = $wsCoverage.Range("R1").Column
$wsCoverage.Range("AE2").Select() | Out-Null

$wb.Save()
